$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.125.06'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '1.788.24'
$ws.Range('E3').Value = '  -1.47%  '
$ws.Range('E4').Value = '  -0.23%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '227.11'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.552'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +1.89%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '31.68'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.31%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '46.28'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.282'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  +0.82%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0660'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -1.31%  '
$ws.Range('E12').Value = '  +0.30%  '
$ws.Range('D13').Value = '2.044.03'
$ws.Range('E13').Value = '  -1.61%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '11.48'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +13.00%  '
$ws.Range('D15').Value = '1.781.30'
$ws.Range('E15').Value = '  -1.87%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.633'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.01%  '
$ws.Range('D17').Value = '34.113.57'
$ws.Range('E17').Value = '  +0.28%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '4.24'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  -1.10%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.52'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.45%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '254.46'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -0.40%  '
$ws.Range('D21').Value = '0.0₃0743'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('E22').Value = '  -0.14%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '10.50'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.98%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '4.25'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.82%  '
$ws.Range('E25').Value = '  -2.03%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '157.22'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -1.45%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '16.60'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '7.04'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -0.52%  '
$ws.Range('E29').Value = '  -0.50%  '
$ws.Range('E30').Value = '  -0.19%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.83'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  -0.24%  '
$ws.Range('E32').Value = '  +1.49%  '
$ws.Range('E33').Value = '  -0.18%  '
$ws.Range('E34').Value = '  +2.19%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.85'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +1.86%  '
$ws.Range('D36').Value = '1.450.31'
$ws.Range('E36').Value = '  -6.86%  '
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.635'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +2.96%  '
$ws.Range('E39').Value = '  +0.45%  '
$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '2.89'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +2.72%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '83.63'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('E42').Value = '  +0.17%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.902'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.09'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -0.91%  '
$ws.Range('E45').Value = '  -2.51%  '
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('B47').Value = 'FraxShare'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '5.88'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +3.40%  '
$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = '1.943.33'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '12.13'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  +9.55%  '
$ws.Range('E50').Value = '  -0.10%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '51.27'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -1.51%  '
